# Update "想去人数" (want-to-go count) values across sheets, matching the
# latest scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 417
$ws1.Range("F4").Value = 2705
$ws1.Range("F10").Value = 581
$ws1.Range("F13").Value = 11240
$ws1.Range("F14").Value = 6414
$ws1.Range("F18").Value = 250
$ws1.Range("F37").Value = 383

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 3644

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8949

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8949
$ws4.Range("F6").Value = 417
$ws4.Range("F7").Value = 2705
$ws4.Range("F16").Value = 581
$ws4.Range("F19").Value = 11240
$ws4.Range("F20").Value = 3644
$ws4.Range("F21").Value = 6414
$ws4.Range("F26").Value = 250
